# Update the "quizvragen" workbook (DC sheet) via Admin:
# append a new multiple-choice question row (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

$ws.Range("B3").Value = "mc"
$ws.Range("D3").Value = "Wat betekend gelijkstroom"
$ws.Range("E3").Value = "['test', 'test1', 'test2']"
$ws.Range("F3").Value = 1
$ws.Range("L3").Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/DC_new_1763470005.jpg"
